$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (AC1) onto the
# three new header cells so they pick up the same bold/border/centered
# style used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every data row (2-43): Wins, Losses, Ties
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 89
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 0
}
